$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Molex Minifit")

# --- Row 27: Molex Micro-Fit 430451601 (16 position) ---
$ws.Cells.Item(27, 1).Value = 430451601
$ws.Cells.Item(27, 2).Value = 16
$ws.Cells.Item(27, 3).Value = "3.0mm"
$ws.Cells.Item(27, 4).Value = "3.0mm"
$ws.Cells.Item(27, 5).Value = 90
$ws.Cells.Item(27, 6).Value = "No"
$ws.Cells.Item(27, 7).Value = "Gold"
$ws.Cells.Item(27, 8).Formula = "=A27"
$ws.Cells.Item(27, 9).Value = "Molex Micro-Fit.SchLib"
$ws.Cells.Item(27, 10).Value = 430451601
$ws.Cells.Item(27, 11).Value = "Molex Micro-Fit.PcbLib"
$ws.Cells.Item(27, 12).Value = "WM7490-ND"

# --- Row 28: Molex Micro-Fit 430451801 (18 position) ---
$ws.Cells.Item(28, 1).Value = 430451801
$ws.Cells.Item(28, 2).Value = 18
$ws.Cells.Item(28, 3).Value = "3.0mm"
$ws.Cells.Item(28, 4).Value = "3.0mm"
$ws.Cells.Item(28, 5).Value = 90
$ws.Cells.Item(28, 6).Value = "No"
$ws.Cells.Item(28, 7).Value = "Gold"
$ws.Cells.Item(28, 8).Value = 430451801
$ws.Cells.Item(28, 9).Value = "Molex Micro-Fit.SchLib"
$ws.Cells.Item(28, 10).Value = 430451801
$ws.Cells.Item(28, 11).Value = "Molex Micro-Fit.PcbLib"
$ws.Cells.Item(28, 12).Value = "2266-430451801-ND"

# Digikey links (M28 before M27, to match the author's original add order)
$ws.Cells.Item(28, 13).Value = "https://www.digikey.com.au/en/products/detail/wec/430451801/18847960"
$ws.Cells.Item(27, 13).Value = "https://www.digikey.com.au/en/products/detail/molex/0430451601/3044582"

$ws.Hyperlinks.Add($ws.Cells.Item(28, 13), "https://www.digikey.com.au/en/products/detail/wec/430451801/18847960") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(27, 13), "https://www.digikey.com.au/en/products/detail/molex/0430451601/3044582") | Out-Null

# Re-apply the same Hyperlink cell style used by the other Digikey-link cells in this column
$ws.Cells.Item(28, 13).Style = "Hyperlink"
$ws.Cells.Item(27, 13).Style = "Hyperlink"

$ws.Select()
$ws.Range("V27").Select()
